$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows at row 18, shifting existing rows down
$ws.Rows.Item(18).Resize(2).Insert()

$ws.Range("A18").Value = "add registry items"
$ws.Range("B18").Value = "Basic"
$ws.Range("A19").Value = "Delete registry items"
$ws.Range("B19").Value = "Basic"

# Selection moves to B20
$ws.Range("B20").Select()

# Print titles: repeat row 1 at top when printing
$ws.PageSetup.PrintTitleRows = "$1:$1"

# Horizontally center on page
$ws.PageSetup.CenterHorizontally = $true
